# Auto-generated edit script: updates LeveProfit-related computed columns (H-N)
# across multiple sheets per the scheduled price-refresh diff.
$wb = $excel.ActiveWorkbook

# === ALC sheet ===
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2221.6875
$ws.Range("J86").Value = 1994.4
$ws.Range("L86").Value = 1994.4
$ws.Range("N86").Value = -4240.4
$ws.Range("H88").Value = 7437.1875
$ws.Range("I88").Value = 1666.6666
$ws.Range("J88").Value = 8768.846
$ws.Range("K88").Value = 1666.6666
$ws.Range("L88").Value = 8768.846
$ws.Range("M88").Value = -1260.6666
$ws.Range("N88").Value = -9580.846
$ws.Range("H89").Value = 2221.6875
$ws.Range("J89").Value = 1994.4
$ws.Range("L89").Value = 9972
$ws.Range("N89").Value = -21204
$ws.Range("H91").Value = 7437.1875
$ws.Range("I91").Value = 1666.6666
$ws.Range("J91").Value = 8768.846
$ws.Range("K91").Value = 1666.6666
$ws.Range("L91").Value = 8768.846
$ws.Range("M91").Value = -262.6666
$ws.Range("N91").Value = -11576.846
$ws.Range("H129").Value = 1096.6154
$ws.Range("J129").Value = 1122.375
$ws.Range("L129").Value = 3367.125
$ws.Range("N129").Value = -13367.125

# === ARM sheet ===
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7764.7334
$ws.Range("I45").Value = 8767
$ws.Range("J45").Value = 1250
$ws.Range("K45").Value = 8767
$ws.Range("L45").Value = 1250
$ws.Range("M45").Value = -8390
$ws.Range("N45").Value = -2004
$ws.Range("H74").Value = 2105.111
$ws.Range("I74").Value = 1986.909
$ws.Range("K74").Value = 1986.909
$ws.Range("M74").Value = -1112.909
$ws.Range("H77").Value = 2105.111
$ws.Range("I77").Value = 1986.909
$ws.Range("K77").Value = 9934.545
$ws.Range("M77").Value = -5566.545
$ws.Range("H102").Value = 5292976.5
$ws.Range("I102").Value = 6174639
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 6174639
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -6173017
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 1787955.9
$ws.Range("I132").Value = 1433.683
$ws.Range("J132").Value = 6671116.5
$ws.Range("K132").Value = 4301.049
$ws.Range("L132").Value = 20013349.5
$ws.Range("M132").Value = -1771.049
$ws.Range("N132").Value = -20018409.5

# === BSM sheet ===
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4402.317
$ws.Range("I134").Value = 4796.7188
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 14390.1564
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -11855.1564
$ws.Range("N134").Value = -14070

# === CRP sheet ===
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4795.871
$ws.Range("I31").Value = 2069.8462
$ws.Range("J31").Value = 6764.6665
$ws.Range("K31").Value = 2069.8462
$ws.Range("L31").Value = 6764.6665
$ws.Range("M31").Value = -1774.8462
$ws.Range("N31").Value = -7354.6665
$ws.Range("H34").Value = 4795.871
$ws.Range("I34").Value = 2069.8462
$ws.Range("J34").Value = 6764.6665
$ws.Range("K34").Value = 2069.8462
$ws.Range("L34").Value = 6764.6665
$ws.Range("M34").Value = -1867.8462
$ws.Range("N34").Value = -7168.6665
$ws.Range("H118").Value = 39789.977
$ws.Range("J118").Value = 39789.977
$ws.Range("L118").Value = 39789.977
$ws.Range("N118").Value = -43103.977
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H130").Value = 52780
$ws.Range("J130").Value = 52780
$ws.Range("L130").Value = 52780
$ws.Range("N130").Value = -62820
$ws.Range("H132").Value = 1601.0714
$ws.Range("I132").Value = 1287.1945
$ws.Range("J132").Value = 3484.3333
$ws.Range("K132").Value = 3861.5835
$ws.Range("L132").Value = 10452.9999
$ws.Range("M132").Value = -1331.5835
$ws.Range("N132").Value = -15512.9999
$ws.Range("H135").Value = 32400
$ws.Range("J135").Value = 32400
$ws.Range("L135").Value = 32400
$ws.Range("N135").Value = -42540
$ws.Range("H138").Value = 65333.332
$ws.Range("J138").Value = 65333.332
$ws.Range("L138").Value = 65333.332
$ws.Range("N138").Value = -75613.33199999999
$ws.Range("H140").Value = 28846.428
$ws.Range("J140").Value = 28846.428
$ws.Range("L140").Value = 28846.428
$ws.Range("N140").Value = -39206.428
$ws.Range("N123").ClearContents()

# === CUL sheet ===
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 264259.75
$ws.Range("I5").Value = 587.2727
$ws.Range("J5").Value = 371681.84
$ws.Range("K5").Value = 1761.8181
$ws.Range("L5").Value = 1115045.52
$ws.Range("M5").Value = -1649.8181
$ws.Range("N5").Value = -1115269.52
$ws.Range("H131").Value = 1667548.1
$ws.Range("J131").Value = 981.3333
$ws.Range("L131").Value = 2943.9999
$ws.Range("N131").Value = -13023.9999
$ws.Range("H135").Value = 264259.75
$ws.Range("I135").Value = 587.2727
$ws.Range("J135").Value = 371681.84
$ws.Range("K135").Value = 5285.454299999999
$ws.Range("L135").Value = 3345136.56
$ws.Range("M135").Value = -2750.454299999999
$ws.Range("N135").Value = -3350206.56

# === GSM sheet ===
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5993333.5
$ws.Range("I7").Value = 17500000
$ws.Range("J7").Value = 240000
$ws.Range("K7").Value = 17500000
$ws.Range("L7").Value = 240000
$ws.Range("M7").Value = -17499888
$ws.Range("N7").Value = -240224
$ws.Range("H8").Value = 5993333.5
$ws.Range("I8").Value = 17500000
$ws.Range("J8").Value = 240000
$ws.Range("K8").Value = 17500000
$ws.Range("L8").Value = 240000
$ws.Range("M8").Value = -17499861
$ws.Range("N8").Value = -240278
$ws.Range("H14").Value = 10000000
$ws.Range("I14").Value = 10000000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 10000000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -9999832
$ws.Range("H126").Value = 6101.077
$ws.Range("I126").Value = 10792
$ws.Range("J126").Value = 2661.0667
$ws.Range("K126").Value = 32376
$ws.Range("L126").Value = 7983.2001
$ws.Range("M126").Value = -29906
$ws.Range("N126").Value = -12923.2001
$ws.Range("H132").Value = 2552.9546
$ws.Range("I132").Value = 1922.1
$ws.Range("K132").Value = 5766.299999999999
$ws.Range("M132").Value = -3236.299999999999
$ws.Range("N14").ClearContents()

# === LTW sheet ===
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 650740.4
$ws.Range("I82").Value = 835501.3
$ws.Range("J82").Value = 207314
$ws.Range("K82").Value = 835501.3
$ws.Range("L82").Value = 207314
$ws.Range("M82").Value = -835140.3
$ws.Range("N82").Value = -208036
$ws.Range("H85").Value = 650740.4
$ws.Range("I85").Value = 835501.3
$ws.Range("J85").Value = 207314
$ws.Range("K85").Value = 835501.3
$ws.Range("L85").Value = 207314
$ws.Range("M85").Value = -834253.3
$ws.Range("N85").Value = -209810
$ws.Range("H132").Value = 14499199
$ws.Range("I132").Value = 22229962
$ws.Range("J132").Value = 4018.375
$ws.Range("K132").Value = 66689886
$ws.Range("L132").Value = 12055.125
$ws.Range("M132").Value = -66687356
$ws.Range("N132").Value = -17115.125
